# Update the "2. Организация" (organization/contact) block on the sheet
# with the new responsible-person details, and re-point the later
# shared-string references (handled automatically by the engine once the
# old, now-unused strings drop out of xl/sharedStrings.xml).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New contact information for rows 6-10 (column B)
$ws.Range("B6").Value = "Национальный статистический комитет КР (Упарвление статистики финансов и цен)."
$ws.Range("B7").Value = "Абдукадирова М.А. `n"
$ws.Range("B8").Value = "m.abdukadirova@stat.kg"
$ws.Range("B9").Value = "(0312) -62-55-91"
$ws.Range("B10").Value = "www.stat.gov.kg"

# The organisation name and contact-person cells now wrap their text
$ws.Range("B6").WrapText = $true
$ws.Range("B7").WrapText = $true

# Move the sheet's active selection to the phone-number cell
[void]$ws.Range("B9").Select()
